$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new status row (row 13) with date 2017-12-08 (serial 43077)
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 43077

$ws.Range("C13").Value = "C# Admin Login: multiple columns in listview with buttons and images (to accept a new entry). Multiple Google Maps (Location of Company) through wrapper class. Both functions without mongoDB (100%)"
$ws.Range("B13").Value = "WebService get collection pupil and entry (90%)."
$ws.Range("D13").Value = "WebService get collection pupil and entry (90%)."

# Update selection to match the new active cell in the diff
$ws.Range("C15").Select()
